$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "0.9999", "241.91") are stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.862.82"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.736.15"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("D5").Value = "241.91"
$ws.Range("E5").Value = "  +4.97%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.5201"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "0.2733"
$ws.Range("D9").Value = "0.06154"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "1.738.79"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "0.07174"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "14.97"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "0.6411"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "4.606"
$ws.Range("D15").Value = "77.10"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "0.9998"
$ws.Range("D18").Value = "25.889.60"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "11.75"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "0.000006774"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "1.961.90"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "4.273"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "8.605"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "5.263"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").Value = "137.36"
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").Value = "1.520"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "15.18"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "1.769"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").Value = "104.93"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").Value = "3.938"
$ws.Range("E30").Value = "  +5.49%  "
$ws.Range("D31").Value = "0.08243"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").Value = "3.653"
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("D33").Value = "0.04649"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").Value = "2.664"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").Value = "0.9885"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "0.6180"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "0.01597"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "1.920"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "0.9998"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "100.35"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "0.3846"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "0.7463"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D44").Value = "4.997"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "0.1123"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "6.247"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "0.05222"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").Value = "54.93"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("D49").Value = "30.55"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").Value = "7.510"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("D51").Value = "0.3408"
$ws.Range("E51").Value = "  -0.65%  "

# Restore the default style on column D so no stray per-cell style
# index is left behind (values remain text because they are already
# stored as shared strings at this point).
$ws.Range("D2:D51").Style = "Normal"
